$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSetup")

# Update test/question codes on row 2 and row 3 (A/B columns)
$ws.Range("A2").Value = "Tery33T"
$ws.Range("B2").Value = "Inregdal3"
$ws.Range("A3").Value = "Terty4T"
$ws.Range("B3").Value = "Inregdal4"

# Move the active selection from C5 to B3
$ws.Range("B3").Select()
